# Trading update: 2026-02-17 04:10:49
# Appends a new "OPEN" MarketMaking trade row (Trade #29, row 30) to both
# the "All Trades" and "MarketMaking" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 30

    # A - Trade # (number)
    $ws.Cells.Item($row, 1).Value = 29

    # B - Date (text). Force text format first so the ISO-like date string
    # ("2026-02-17") isn't auto-converted into a date serial number.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"

    # C - Time (text)
    $ws.Cells.Item($row, 3).Value = "04:09:11"

    # D - Strategy (text)
    $ws.Cells.Item($row, 4).Value = "MarketMaking"

    # E - Side (text)
    $ws.Cells.Item($row, 5).Value = "DOWN"

    # F - Entry Price (number)
    $ws.Cells.Item($row, 6).Value = 0.53

    # G - Exit Price -> blank (trade still OPEN, no exit price yet)

    # H - Status (text)
    $ws.Cells.Item($row, 8).Value = "OPEN"

    # I - P&L % (number)
    $ws.Cells.Item($row, 9).Value = 0

    # J - P&L $ (number)
    $ws.Cells.Item($row, 10).Value = 0

    # K - Capital After (number)
    $ws.Cells.Item($row, 11).Value = 100.4254564381429

    # L - Entry Slippage (bps) (number)
    $ws.Cells.Item($row, 12).Value = 0

    # M - Exit Slippage (bps) (number)
    $ws.Cells.Item($row, 13).Value = 0

    # N - Confidence (number)
    $ws.Cells.Item($row, 14).Value = 0.6

    # O - Entry Reason (text)
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    # P - Exit Reason -> blank (trade still OPEN, no exit reason yet)

    # Q - Duration (min) (number)
    $ws.Cells.Item($row, 17).Value = 0
}
